$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both contain the same two data rows that need updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 472
    $ws.Range("F3").Value = 45
}
